$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add "File Name" in A1 (matching the style used by B1/C1),
# and rename the B1/C1 headers.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Value = "File Name"

$ws.Range("B1").Value = "Toplam Katma Değer Vergisi"
$ws.Range("C1").Value = "Matrah Toplamı"

# Row 2: move the file name into column A (no special style) and shift
# the numeric/value columns over.
$ws.Range("A2").Value = "Örnek KDV Beyannamesi II.pdf"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "113.122.268,06"
$ws.Range("C2").Value = "2.936.636,47"

# Row 3: move the file name into column A (no special style) and shift
# the numeric/value columns over.
$ws.Range("A3").Value = "Örnek KDV Beyannamesi.pdf"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "565.316.718,40"
$ws.Range("C3").Value = "14.683.182,35"
